$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.357.66"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.30%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.650.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -3.83%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.510"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.93"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("E9").Value = "  -1.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0615"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.886.66"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.651.66"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.562"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.401.01"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "234.71"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.50"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.22"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.92%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.03"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("E31").Value = "  -0.96%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.30"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.447.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.06"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.909"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.569"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.27%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.58"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.43"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("E44").Value = "  -2.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.789"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.81%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.794.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.27"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0101"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -11.44%  "
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.77"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.53%  "
